$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to be
# bumped from 45177 (2023-09-08) to 45178 (2023-09-09) for every data
# row (rows 2 through 358), leaving the cell style/format untouched.
$ws.Range("C2:C358").Value = 45178
